# Apply the "Trade #32 closed" update described by the commit:
#   Trade #32 closed at 2026-02-17 08:04:07 - unknown UNKNOWN +0.000%
#
# This touches four sheets:
#   Summary         -> roll the aggregate stats (Current Capital, Total P&L $/%,
#                      Total Trades, Losing Trades, Win Rate %)
#   Strategy Status -> roll the MarketMaking strategy row (Capital, Trades,
#                      P&L $, P&L %, Win Rate %)
#   All Trades      -> append the new trade as row 33
#   MarketMaking    -> append the same trade as row 33 (per-strategy trade log)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a text value into a cell without Excel's implicit
# string->date coercion (e.g. "2026-02-17" becoming a date serial number)
# and without leaving a residual NumberFormat/style behind on the cell.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.37   # Current Capital
$wsSummary.Range("B4").Value = -0.63     # Total P&L $
$wsSummary.Range("B5").Value = -0.39     # Total P&L %
$wsSummary.Range("B6").Value = 32        # Total Trades
$wsSummary.Range("B8").Value = 16        # Losing Trades
$wsSummary.Range("B9").Value = 25        # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row, row 4)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.37      # Capital
$wsStatus.Range("D4").Value = 32         # Trades
$wsStatus.Range("E4").Value = -0.63      # P&L $
$wsStatus.Range("F4").Value = -0.63      # P&L %
$wsStatus.Range("G4").Value = 25         # Win Rate %

# ---------------------------------------------------------------------------
# New trade record (trade #32) appended to both "All Trades" and
# "MarketMaking" sheets as row 33.
# ---------------------------------------------------------------------------
$sheetsToAppend = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetsToAppend) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A33").Value = 32

    # "2026-02-17" looks like a date to Excel's Range.Value setter and would
    # otherwise be silently coerced into a date serial number - force it to
    # stay text. The other text columns below (time-of-day, strategy name,
    # side, status, reasons) are not subject to that coercion so they can be
    # assigned directly.
    Set-TextValue $ws.Range("B33") "2026-02-17"
    $ws.Range("C33").Value = "08:04:01"
    $ws.Range("D33").Value = "MarketMaking"
    $ws.Range("E33").Value = "UP"

    $ws.Range("F33").Value = 0.27
    $ws.Range("G33").Value = 0.26

    $ws.Range("H33").Value = "CLOSED"

    $ws.Range("I33").Value = -3.7037
    $ws.Range("J33").Value = -0.01
    $ws.Range("K33").Value = 99.37
    $ws.Range("L33").Value = 0
    $ws.Range("M33").Value = 0
    $ws.Range("N33").Value = 0.6

    $ws.Range("O33").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P33").Value = "early_exit"

    $ws.Range("Q33").Value = 0.13
}

Write-Output "applied trade #32 close update"
